# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets
# to reflect freshly generated numbers (commit: "Update gh-pages to output
# generated at 456a3b4").

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1437
$ws1.Range("F3").Value = 3023
$ws1.Range("F4").Value = 35
$ws1.Range("F5").Value = 298
$ws1.Range("F6").Value = 283

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1437
$ws4.Range("F3").Value = 3023
$ws4.Range("F4").Value = 35
$ws4.Range("F5").Value = 298
$ws4.Range("F7").Value = 283
